$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 (sldId 259 / cId 2797378828) - Content Placeholder 2 (Shape id=3)
# Merge "Data Source:" + "ETL" paragraphs into a single "Data Source and ETL"
# paragraph.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange

$dsPara = $tr2.Paragraphs(3, 1)
$etlPara = $tr2.Paragraphs(4, 1)
$mergedRange = $tr2.Characters($dsPara.Start, ($etlPara.Start + $etlPara.Length - $dsPara.Start))
$mergedRange.Text = "Data Source and ETL"
# Remove the now-empty paragraph left behind by the text merge above.
$tr2.Paragraphs(4, 1).Delete()

# ---------------------------------------------------------------------------
# Slide 4 (sldId 261 / cId 481207256) - Content Placeholder 2 (Shape id=3)
# "Library: " + "SKLearn" -> "Library: SKLearn"
# Collapse the long algorithm list into a single run (after "LogisticRegression").
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

$libPara = $tr4.Paragraphs(4, 1)
$libRange = $tr4.Characters($libPara.Start, $libPara.Length - 1)
$libRange.Text = "Library: SKLearn"

$algoPara = $tr4.Paragraphs(5, 1)
$afterLogisticStart = $algoPara.Start + 18
$afterLogisticLen = $algoPara.Length - 18 - 1
$restRange = $tr4.Characters($afterLogisticStart, $afterLogisticLen)
$restRange.Text = ", DecisionTree, RandomForest, GradientBoosting, SGD,  SVM, MultiNomialNB, Metrics"

# ---------------------------------------------------------------------------
# Slide 5 (sldId 257 / cId 2114131046) - Content Placeholder 2 (Shape id=3)
# "Hover filter" -> "Click filter"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange

$funcPara = $tr5.Paragraphs(5, 1)
$funcRange = $tr5.Characters($funcPara.Start, $funcPara.Length - 1)
$funcRange.Text = "Tableau Functions Used: Click filter, Cluster Model (Analytics), LOD"

# ---------------------------------------------------------------------------
# Slide 6 (sldId 258 / cId 3374226764) - Content Placeholder 2 (Shape id=3)
# "Link to Dashboards" + ": " -> "Link to Dashboards: " (merge into one run,
# leave the hyperlink run untouched).
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6 = $shp6.TextFrame.TextRange

$linkPara = $tr6.Paragraphs(5, 1)
$prefixRange = $tr6.Characters($linkPara.Start, 20)
$prefixRange.Text = "Link to Dashboards: "
